# Applies the cryptos list refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.847.30'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '''1.894.16'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = '''0.7878'
$ws.Range("E5").Value = '  -6.40%  '
$ws.Range("D6").Value = '''243.35'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").Value = '''0.3143'
$ws.Range("E8").Value = '  -4.66%  '
$ws.Range("D9").Value = '''25.26'
$ws.Range("E9").Value = '  -5.67%  '
$ws.Range("D10").Value = '''0.07255'
$ws.Range("E10").Value = '  +2.35%  '
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("D12").Value = '''0.7644'
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").Value = '''5.534'
$ws.Range("E13").Value = '  +4.66%  '
$ws.Range("D14").Value = '''1.905.84'
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").Value = '''92.38'
$ws.Range("E15").Value = '  -0.46%  '
$ws.Range("D16").Value = '''6.132'
$ws.Range("E16").Value = '  +4.02%  '
$ws.Range("D17").Value = '''29.862.98'
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("E18").Value = '  -2.09%  '
$ws.Range("D19").Value = '''244.08'
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("D20").Value = '''0.000007780'
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '''1.001'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '''2.153.74'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").Value = '''8.117'
$ws.Range("E23").Value = '  +15.51%  '
$ws.Range("D24").Value = '''1.002'
$ws.Range("E24").Value = '  +0.32%  '
$ws.Range("D25").Value = '''0.1645'
$ws.Range("E25").Value = '  -8.25%  '
$ws.Range("D26").Value = '''9.382'
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("D27").Value = '''163.43'
$ws.Range("E27").Value = '  -1.29%  '
$ws.Range("D28").Value = '''18.71'
$ws.Range("E28").Value = '  -1.58%  '
$ws.Range("D29").Value = '''2.050'
$ws.Range("E29").Value = '  -2.84%  '
$ws.Range("D30").Value = '''1.400'
$ws.Range("E30").Value = '  +2.47%  '
$ws.Range("D31").Value = '''1.546'
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("D32").Value = '''4.465'
$ws.Range("E32").Value = '  +3.69%  '
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("D34").Value = '''0.05538'
$ws.Range("E34").Value = '  -6.99%  '
$ws.Range("D35").Value = '''1.267'
$ws.Range("E35").Value = '  -0.62%  '
$ws.Range("D36").Value = '''0.7383'
$ws.Range("E36").Value = '  +0.55%  '
$ws.Range("D37").Value = '''0.9992'
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("D38").Value = '''2.631'
$ws.Range("E38").Value = '  -2.97%  '
$ws.Range("D39").Value = '''0.01924'
$ws.Range("E39").Value = '  -0.33%  '
$ws.Range("D40").Value = '''2.780'
$ws.Range("E40").Value = '  -0.27%  '
$ws.Range("D41").Value = '''1.141.78'
$ws.Range("E41").Value = '  +13.49%  '
$ws.Range("D42").Value = '''73.93'
$ws.Range("E42").Value = '  +1.05%  '
$ws.Range("D43").Value = '''0.4419'
$ws.Range("E43").Value = '  -0.93%  '
$ws.Range("D44").Value = '''5.880'
$ws.Range("E44").Value = '  -1.36%  '
$ws.Range("D45").Value = '''0.8506'
$ws.Range("E45").Value = '  -0.81%  '
$ws.Range("D46").Value = '''104.28'
$ws.Range("E46").Value = '  +1.99%  '
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''10.04'
$ws.Range("E48").Value = '  +1.87%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''1.876'
$ws.Range("E49").Value = '  -1.88%  '
$ws.Range("D50").Value = '''7.440'
$ws.Range("E50").Value = '  -2.04%  '
$ws.Range("D51").Value = '''3.002'
$ws.Range("E51").Value = '  +9.68%  '
